$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country name labels (re-sorted rows) ---
$ws.Range("A41").Value = "Republica Dominicana"
$ws.Range("A42").Value = "Israel"

$ws.Range("A81").Value = "Republica de Macedonia"
$ws.Range("A82").Value = "Bulgaria"

$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"

# --- Update case-count figures for affected countries ---
$ws.Range("B4").Value = 4122702
$ws.Range("C4").Value = 21827
$ws.Range("D4").Value = 1944490
$ws.Range("E4").Value = 2031664
$ws.Range("G4").Value = 365
$ws.Range("H4").Value = 146548
$ws.Range("B5").Value = 2242394
$ws.Range("C5").Value = 10523
$ws.Range("E5").Value = 627220
$ws.Range("G5").Value = 146
$ws.Range("H5").Value = 83036
$ws.Range("B6").Value = 1284638
$ws.Range("C6").Value = 44954
$ws.Range("D6").Value = 814912
$ws.Range("E6").Value = 439125
$ws.Range("G6").Value = 711
$ws.Range("H6").Value = 30601
$ws.Range("B12").Value = 317246
$ws.Range("C12").Value = 2615
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 28429
$ws.Range("B13").Value = 297146
$ws.Range("C13").Value = 769
$ws.Range("G13").Value = 53
$ws.Range("H13").Value = 45554
$ws.Range("B17").Value = 245338
$ws.Range("C17").Value = 306
$ws.Range("D17").Value = 197842
$ws.Range("E17").Value = 12404
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 35092
$ws.Range("B25").Value = 108244
$ws.Range("C25").Value = 373
$ws.Range("D25").Value = 105018
$ws.Range("E25").Value = 3062
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 164
$ws.Range("B41").Value = 57615
$ws.Range("C41").Value = 1572
$ws.Range("D41").Value = 26905
$ws.Range("E41").Value = 29704
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = 1006
$ws.Range("B42").Value = 56748
$ws.Range("C42").Value = 663
$ws.Range("D42").Value = 23560
$ws.Range("E42").Value = 32755
$ws.Range("G42").Value = 3
$ws.Range("H42").Value = 433
$ws.Range("B71").Value = 14724
$ws.Range("C71").Value = 154
$ws.Range("D71").Value = 9292
$ws.Range("E71").Value = 5067
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 365
$ws.Range("B78").Value = 11933
$ws.Range("C78").Value = 409
$ws.Range("D78").Value = 5645
$ws.Range("E78").Value = 6091
$ws.Range("G78").Value = 9
$ws.Range("H78").Value = 197
$ws.Range("B81").Value = 9669
$ws.Range("C81").Value = 122
$ws.Range("D81").Value = 5071
$ws.Range("E81").Value = 4153
$ws.Range("G81").Value = 3
$ws.Range("H81").Value = 445
$ws.Range("B82").Value = 9584
$ws.Range("D82").Value = 4643
$ws.Range("E82").Value = 4620
$ws.Range("H82").Value = 321
$ws.Range("B96").Value = 5952
$ws.Range("C96").Value = 98
$ws.Range("D96").Value = 4591
$ws.Range("E96").Value = 1249
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 112
$ws.Range("B102").Value = 4110
$ws.Range("C102").Value = 33
$ws.Range("E102").Value = 2535
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 201
$ws.Range("B130").Value = 1752
$ws.Range("C130").Value = 21
$ws.Range("D130").Value = 1292
$ws.Range("E130").Value = 394
$ws.Range("B134").Value = 1582
$ws.Range("C134").Value = 25
$ws.Range("D134").Value = 528
$ws.Range("E134").Value = 1043
$ws.Range("B136").Value = 1522
$ws.Range("C136").Value = 120
$ws.Range("D136").Value = 69
$ws.Range("E136").Value = 1446
$ws.Range("B137").Value = 1406
$ws.Range("C137").Value = 12
$ws.Range("D137").Value = 1118
$ws.Range("E137").Value = 238
$ws.Range("B158").Value = 584
$ws.Range("C158").Value = 23
$ws.Range("D158").Value = 174
$ws.Range("E158").Value = 375
$ws.Range("G158").Value = 3
$ws.Range("H158").Value = 35

# --- Update "last updated" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Julio de 2020 a las 18:30"
